$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "70.136.54"
$ws.Range("E2").Value = "  -0.99%  "
# Row 3
$ws.Range("D3").Value = "3.577.31"
$ws.Range("E3").Value = "  -1.71%  "
# Row 4
$ws.Range("E4").Value = "  +0.10%  "
# Row 5
Set-TextValue $ws.Range("D5") "578.74"
$ws.Range("E5").Value = "  -2.51%  "
# Row 6
Set-TextValue $ws.Range("D6") "186.51"
$ws.Range("E6").Value = "  -4.44%  "
# Row 7
$ws.Range("D7").Value = "3.573.14"
$ws.Range("E7").Value = "  -1.66%  "
# Row 8
Set-TextValue $ws.Range("D8") "0.619"
$ws.Range("E8").Value = "  -4.11%  "
# Row 9
$ws.Range("E9").Value = "  +0.11%  "
# Row 10
Set-TextValue $ws.Range("D10") "0.183"
$ws.Range("E10").Value = "  -1.43%  "
# Row 11
Set-TextValue $ws.Range("D11") "0.651"
$ws.Range("E11").Value = "  -4.48%  "
# Row 12
Set-TextValue $ws.Range("D12") "55.10"
$ws.Range("E12").Value = "  -5.22%  "
# Row 13
Set-TextValue $ws.Range("D13") "0.0000305"
$ws.Range("E13").Value = "  -1.14%  "
# Row 14
Set-TextValue $ws.Range("D14") "9.54"
$ws.Range("E14").Value = "  -4.51%  "
# Row 15
$ws.Range("D15").Value = "4.152.31"
$ws.Range("E15").Value = "  -1.67%  "
# Row 16
Set-TextValue $ws.Range("D16") "19.66"
$ws.Range("E16").Value = "  -4.05%  "
# Row 17
$ws.Range("D17").Value = "3.581.54"
$ws.Range("E17").Value = "  -1.56%  "
# Row 18
$ws.Range("D18").Value = "70.075.49"
$ws.Range("E18").Value = "  -1.01%  "
# Row 19
Set-TextValue $ws.Range("D19") "12.57"
$ws.Range("E19").Value = "  -1.76%  "
# Row 20
Set-TextValue $ws.Range("D20") "0.120"
$ws.Range("E20").Value = "  -1.18%  "
# Row 21
$ws.Range("E21").Value = "  -3.22%  "
# Row 22
Set-TextValue $ws.Range("D22") "493.38"
$ws.Range("E22").Value = "  +0.92%  "
# Row 23
Set-TextValue $ws.Range("D23") "19.32"
$ws.Range("E23").Value = "  -0.33%  "
# Row 24
$ws.Range("E24").Value = "  -5.34%  "
# Row 25
Set-TextValue $ws.Range("D25") "96.76"
$ws.Range("E25").Value = "  +5.77%  "
# Row 26
$ws.Range("E26").Value = "  -2.42%  "
# Row 27
Set-TextValue $ws.Range("D27") "11.49"
$ws.Range("E27").Value = "  +0.50%  "
# Row 28
$ws.Range("E28").Value = "  -6.97%  "
# Row 29
Set-TextValue $ws.Range("D29") "9.34"
$ws.Range("E29").Value = "  -2.71%  "
# Row 30
Set-TextValue $ws.Range("D30") "7.77"
$ws.Range("E30").Value = "  -2.30%  "
# Row 31
Set-TextValue $ws.Range("D31") "31.69"
$ws.Range("E31").Value = "  -3.51%  "
# Row 32
Set-TextValue $ws.Range("D32") "12.12"
$ws.Range("E32").Value = "  -1.51%  "
# Row 33
Set-TextValue $ws.Range("D33") "65.92"
$ws.Range("E33").Value = "  -0.74%  "
# Row 34
Set-TextValue $ws.Range("D34") "0.115"
$ws.Range("E34").Value = "  -6.10%  "
# Row 35
Set-TextValue $ws.Range("D35") "580.90"
$ws.Range("E35").Value = "  -5.46%  "
# Row 36
Set-TextValue $ws.Range("D36") "3.22"
$ws.Range("E36").Value = "  +14.29%  "
# Row 39
$ws.Range("E39").Value = "  +0.10%  "
# Row 40
$ws.Range("D40").Value = "0.0₃0794"
$ws.Range("E40").Value = "  -4.69%  "
# Row 41
Set-TextValue $ws.Range("D41") "3.48"
$ws.Range("E41").Value = "  -2.57%  "
# Row 42
$ws.Range("E42").Value = "  -0.01%  "
# Row 43
$ws.Range("E43").Value = "  -9.25%  "
# Row 44
Set-TextValue $ws.Range("D44") "3.58"
$ws.Range("E44").Value = "  +7.39%  "
# Row 45
Set-TextValue $ws.Range("D45") "3.06"
$ws.Range("E45").Value = "  -4.02%  "
# Row 46
Set-TextValue $ws.Range("D46") "0.0442"
$ws.Range("E46").Value = "  -3.81%  "
# Row 47
$ws.Range("D47").Value = "3.185.94"
$ws.Range("E47").Value = "  -4.34%  "
# Row 48
$ws.Range("E48").Value = "  -2.09%  "
# Row 49
Set-TextValue $ws.Range("D49") "1.58"
$ws.Range("E49").Value = "  +31.62%  "
# Row 50
$ws.Range("E50").Value = "  -2.52%  "
# Row 51
Set-TextValue $ws.Range("D51") "1.00"
$ws.Range("E51").Value = "  +0.14%  "
# Row 37 (content swapped with row 38)
$ws.Range("B37").Value = "InjectiveProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D37") "38.93"
$ws.Range("E37").Value = "  -3.54%  "
# Row 38 (content swapped with row 37)
$ws.Range("B38").Value = "TheGraph"
$ws.Range("C38").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue $ws.Range("D38") "0.412"
$ws.Range("E38").Value = "  +0.04%  "
